# Apply crypto price/volume updates as per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "43.373.98"
Set-TextValue $ws.Cells.Item(2, 5) "  -1.53%  "

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "2.344.70"
Set-TextValue $ws.Cells.Item(3, 5) "  +3.53%  "

# Row 4
Set-TextValue $ws.Cells.Item(4, 5) "  -0.03%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "0.652"
Set-TextValue $ws.Cells.Item(5, 5) "  +1.97%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "231.13"
Set-TextValue $ws.Cells.Item(6, 5) "  +0.12%  "

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "65.44"
Set-TextValue $ws.Cells.Item(7, 5) "  +1.24%  "

# Row 8
Set-TextValue $ws.Cells.Item(8, 5) "  -0.04%  "

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "0.459"
Set-TextValue $ws.Cells.Item(9, 5) "  +2.35%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "0.0952"
Set-TextValue $ws.Cells.Item(10, 5) "  -4.83%  "

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "56.90"
Set-TextValue $ws.Cells.Item(11, 5) "  -0.66%  "

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "26.72"
Set-TextValue $ws.Cells.Item(12, 5) "  -2.25%  "

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "2.689.43"
Set-TextValue $ws.Cells.Item(13, 5) "  +3.36%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 5) "  -2.13%  "

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "15.35"
Set-TextValue $ws.Cells.Item(15, 5) "  -2.51%  "

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "6.26"
Set-TextValue $ws.Cells.Item(16, 5) "  +2.70%  "

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) "0.842"
Set-TextValue $ws.Cells.Item(17, 5) "  +0.20%  "

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "2.339.62"
Set-TextValue $ws.Cells.Item(18, 5) "  +3.07%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "43.280.37"
Set-TextValue $ws.Cells.Item(19, 5) "  -1.52%  "

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) "0.0₃0977"
Set-TextValue $ws.Cells.Item(20, 5) "  -3.34%  "

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "73.61"
Set-TextValue $ws.Cells.Item(21, 5) "  -0.20%  "

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "6.18"
Set-TextValue $ws.Cells.Item(22, 5) "  +1.01%  "

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) "248.21"
Set-TextValue $ws.Cells.Item(23, 5) "  -1.89%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "3.91"
Set-TextValue $ws.Cells.Item(24, 5) "  +20.55%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 5) "  +0.02%  "

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) "2.44"
Set-TextValue $ws.Cells.Item(26, 5) "  -0.88%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 5) "  +0.79%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) "9.91"
Set-TextValue $ws.Cells.Item(28, 5) "  -2.07%  "

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) "175.88"
Set-TextValue $ws.Cells.Item(29, 5) "  +2.52%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "22.22"
Set-TextValue $ws.Cells.Item(30, 5) "  +5.98%  "

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) "1.51"
Set-TextValue $ws.Cells.Item(31, 5) "  +7.89%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "0.129"
Set-TextValue $ws.Cells.Item(32, 5) "  -8.07%  "

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) "0.126"
Set-TextValue $ws.Cells.Item(33, 5) "  +0.14%  "

# Row 34
Set-TextValue $ws.Cells.Item(34, 5) "  +4.00%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "0.0686"
Set-TextValue $ws.Cells.Item(35, 5) "  -3.07%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) "4.99"
Set-TextValue $ws.Cells.Item(36, 5) "  +1.41%  "

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) "2.48"
Set-TextValue $ws.Cells.Item(37, 5) "  +6.95%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) "6.47"
Set-TextValue $ws.Cells.Item(38, 5) "  -1.04%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "3.59"
Set-TextValue $ws.Cells.Item(39, 5) "  -6.13%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "0.0252"
Set-TextValue $ws.Cells.Item(40, 5) "  -3.34%  "

# Row 41
Set-TextValue $ws.Cells.Item(41, 5) "  +0.09%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "8.89"
Set-TextValue $ws.Cells.Item(42, 5) "  +7.93%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "17.93"
Set-TextValue $ws.Cells.Item(43, 5) "  +2.32%  "

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "1.16"
Set-TextValue $ws.Cells.Item(44, 5) "  +6.93%  "

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) "98.66"
Set-TextValue $ws.Cells.Item(45, 5) "  +0.18%  "

# Row 46
Set-TextValue $ws.Cells.Item(46, 5) "  -0.83%  "

# Row 47
Set-TextValue $ws.Cells.Item(47, 2) "Cronos"
Set-TextValue $ws.Cells.Item(47, 3) "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Cells.Item(47, 4) "0.0946"
Set-TextValue $ws.Cells.Item(47, 5) "  -3.68%  "

# Row 48
Set-TextValue $ws.Cells.Item(48, 2) "FTXToken"
Set-TextValue $ws.Cells.Item(48, 3) "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Cells.Item(48, 4) "4.38"
Set-TextValue $ws.Cells.Item(48, 5) "  -1.42%  "

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "1.436.59"
Set-TextValue $ws.Cells.Item(49, 5) "  -0.80%  "

# Row 50
Set-TextValue $ws.Cells.Item(50, 2) "TerraClassic"
Set-TextValue $ws.Cells.Item(50, 3) "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
Set-TextValue $ws.Cells.Item(50, 4) "0.000203"
Set-TextValue $ws.Cells.Item(50, 5) "  -8.29%  "

# Row 51
Set-TextValue $ws.Cells.Item(51, 2) "Celestia"
Set-TextValue $ws.Cells.Item(51, 3) "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Cells.Item(51, 4) "9.83"
Set-TextValue $ws.Cells.Item(51, 5) "  -6.31%  "
